# Update "Generate Report for Handback" timestamps on the zh-cn and de-de
# sheets for the 8e2dfd0b-... row (row 2).
#
# Correspond Handoff Datetime (column D) and Correspond Handback DateTime
# (column G) are refreshed to reflect the new report generation time.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-22 16:42:29"
$wsZhCn.Range("G2").Value = "2016-02-22 16:43:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-22 16:42:39"
$wsDeDe.Range("G2").Value = "2016-02-22 16:43:46"
